$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of repo/file/search/replace data appended to the table
$ws.Range("A9").Value = "github-gk-aks/thirdgithubrepo"
$ws.Range("B9").Value = ".github/workflows/github-rest-api1234.yaml"
$ws.Range("C9").Value = "A_TOKEN"
$ws.Range("D9").Value = "B_TOKEN"

$ws.Range("A10").Value = "github-gk-aks/thirdgithubrepo"
$ws.Range("B10").Value = ".github/workflows/github-rest-api123.yaml"
$ws.Range("C10").Value = "A_TOKEN"
$ws.Range("D10").Value = "B_TOKEN"

$ws.Range("A11").Value = "github-gk-aks/fourthgithubrepo"
$ws.Range("B11").Value = ".github/workflows/github-rest-api1.yaml"
$ws.Range("C11").Value = "A_TOKEN"
$ws.Range("D11").Value = "C_TOKEN"

$ws.Range("A12").Value = "github-gk-aks/fourthgithubrepo"
$ws.Range("B12").Value = ".github/workflows/github-rest-api1.yaml"
$ws.Range("C12").Value = "A_TOKEN"
$ws.Range("D12").Value = "C_TOKEN"

$ws.Range("A13").Value = "github-gk-aks/fifthgithubrepo"
$ws.Range("B13").Value = ".github/workflows/github-rest-api1234.yaml"
$ws.Range("C13").Value = "A_TOKEN"
$ws.Range("D13").Value = "D_TOKEN"

# Match the final selection state from the saved workbook
$ws.Range("C12:D13").Select()
